$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "58.968.39"
Set-TextValue $ws.Range("E2") "  +2.91%  "
Set-TextValue $ws.Range("D3") "2.587.31"
Set-TextValue $ws.Range("E3") "  +1.29%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "521.16"
Set-TextValue $ws.Range("E5") "  +1.04%  "
Set-TextValue $ws.Range("D6") "139.78"
Set-TextValue $ws.Range("E6") "  -1.40%  "
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("D8") "0.564"
Set-TextValue $ws.Range("E8") "  +0.10%  "
Set-TextValue $ws.Range("D9") "2.600.35"
Set-TextValue $ws.Range("E9") "  +1.38%  "
Set-TextValue $ws.Range("E10") "  -0.54%  "
Set-TextValue $ws.Range("E11") "  +1.25%  "
Set-TextValue $ws.Range("E12") "  +2.66%  "
Set-TextValue $ws.Range("E13") "  +3.25%  "
Set-TextValue $ws.Range("D14") "3.045.00"
Set-TextValue $ws.Range("E14") "  +1.29%  "
Set-TextValue $ws.Range("D15") "58.981.45"
Set-TextValue $ws.Range("E15") "  +2.96%  "
Set-TextValue $ws.Range("D16") "20.46"
Set-TextValue $ws.Range("E16") "  +2.03%  "
Set-TextValue $ws.Range("B17") "ShibaInu"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D17") "0.0000133"
Set-TextValue $ws.Range("E17") "  +0.48%  "
Set-TextValue $ws.Range("B18") "WrappedEther"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D18") "2.581.83"
Set-TextValue $ws.Range("E18") "  +1.64%  "
Set-TextValue $ws.Range("D19") "339.05"
Set-TextValue $ws.Range("E19") "  +2.02%  "
Set-TextValue $ws.Range("D20") "4.31"
Set-TextValue $ws.Range("E20") "  +1.31%  "
Set-TextValue $ws.Range("D21") "10.11"
Set-TextValue $ws.Range("E21") "  +0.12%  "
Set-TextValue $ws.Range("D22") "6.47"
Set-TextValue $ws.Range("E22") "  +3.42%  "
Set-TextValue $ws.Range("D23") "0.999"
Set-TextValue $ws.Range("D24") "66.14"
Set-TextValue $ws.Range("E24") "  +1.76%  "
Set-TextValue $ws.Range("D25") "0.167"
Set-TextValue $ws.Range("E25") "  +1.21%  "
Set-TextValue $ws.Range("E26") "  +1.20%  "
Set-TextValue $ws.Range("D27") "0.994"
Set-TextValue $ws.Range("E27") "  -0.38%  "
Set-TextValue $ws.Range("E28") "  +1.89%  "
Set-TextValue $ws.Range("E29") "  +0.00%  "
Set-TextValue $ws.Range("D30") "0.0₃0727"
Set-TextValue $ws.Range("E30") "  -1.84%  "
Set-TextValue $ws.Range("D31") "5.94"
Set-TextValue $ws.Range("E31") "  -5.96%  "
Set-TextValue $ws.Range("E32") "  +0.33%  "
Set-TextValue $ws.Range("D33") "18.72"
Set-TextValue $ws.Range("E33") "  +1.10%  "
Set-TextValue $ws.Range("D34") "149.16"
Set-TextValue $ws.Range("E34") "  -0.05%  "
Set-TextValue $ws.Range("D35") "3.99"
Set-TextValue $ws.Range("E35") "  +0.72%  "
Set-TextValue $ws.Range("E36") "  -0.76%  "
Set-TextValue $ws.Range("D37") "36.77"
Set-TextValue $ws.Range("E37") "  +2.46%  "
Set-TextValue $ws.Range("D38") "1.46"
Set-TextValue $ws.Range("E38") "  +2.10%  "
Set-TextValue $ws.Range("D39") "0.827"
Set-TextValue $ws.Range("E39") "  +0.22%  "
Set-TextValue $ws.Range("D40") "0.813"
Set-TextValue $ws.Range("E40") "  -5.54%  "
Set-TextValue $ws.Range("D41") "3.50"
Set-TextValue $ws.Range("E41") "  +0.15%  "
Set-TextValue $ws.Range("E42") "  -0.03%  "
Set-TextValue $ws.Range("D43") "273.87"
Set-TextValue $ws.Range("E43") "  +1.80%  "
Set-TextValue $ws.Range("E44") "  +0.92%  "
Set-TextValue $ws.Range("D45") "0.591"
Set-TextValue $ws.Range("E45") "  +1.72%  "
Set-TextValue $ws.Range("D46") "0.0952"
Set-TextValue $ws.Range("E46") "  +0.14%  "
Set-TextValue $ws.Range("D47") "0.0517"
Set-TextValue $ws.Range("E47") "  -0.20%  "
Set-TextValue $ws.Range("D48") "18.48"
Set-TextValue $ws.Range("E48") "  -0.86%  "
Set-TextValue $ws.Range("D49") "1.970.31"
Set-TextValue $ws.Range("E49") "  +0.68%  "
Set-TextValue $ws.Range("D50") "4.54"
Set-TextValue $ws.Range("E50") "  +1.42%  "
Set-TextValue $ws.Range("E51") "  +0.77%  "
